# New weekly record for "Feria Lagunitas de Puerto Montt - Mango":
# Insert a new data row right after the existing row 165 (so it becomes the
# new row 165), pushing every following row down by one. This grows the
# sheet from 214 to 215 used rows (A1:T214 -> A1:T215).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 165 - everything from old row 165 downward
# (including its formatting) shifts down to row 166 onward.
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new weekly record.
$ws.Cells.Item(165, 1).Value = 4
$ws.Cells.Item(165, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(165, 3).Value = "Los Lagos"
$ws.Cells.Item(165, 4).Value = 44754
$ws.Cells.Item(165, 5).Value = 10
$ws.Cells.Item(165, 6).Value = "Fruta"
$ws.Cells.Item(165, 7).Value = 100108
$ws.Cells.Item(165, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(165, 9).Value = 100108002
$ws.Cells.Item(165, 10).Value = "Mango"
$ws.Cells.Item(165, 11).Value = "Sin especificar"
$ws.Cells.Item(165, 12).Value = "Primera"
$ws.Cells.Item(165, 13).Value = 120
$ws.Cells.Item(165, 14).Value = 9000
$ws.Cells.Item(165, 15).Value = 9000
$ws.Cells.Item(165, 16).Value = 9000
$ws.Cells.Item(165, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(165, 18).Value = "Perú"
$ws.Cells.Item(165, 19).Value = 2250
$ws.Cells.Item(165, 20).Value = 4

# Make sure the date cell keeps the same date-style formatting as the
# other rows in column D (the Insert() above already copies row 165's
# old formatting down, so the newly-blank row 165 should inherit it too,
# but set it explicitly to be safe).
$ws.Cells.Item(165, 4).NumberFormat = $ws.Cells.Item(166, 4).NumberFormat()
